$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("network_optimized_weights")
$ws.Activate()

# New edge-name labels, in left-to-right / top-to-bottom order (B1..O1 and A2..A15)
$labels = @("E14", "E1", "E2", "E3", "E4", "E5", "E6", "E7", "E8", "E9", "E10", "E11", "E12", "E13")

for ($i = 0; $i -lt $labels.Length; $i++) {
    # Header row: columns B..O (column index 2..15)
    $ws.Cells.Item(1, $i + 2).Value = $labels[$i]
    # Row labels: rows 2..15, column A
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}

# Update the active selection shown in the sheet view
$ws.Range("C20").Select()
